# Ausimplementierung Produktvorschlag + Update Embedding
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Laufzeit (E) values first: mittelfristig (rows 2 & 3), then kurzfristig (row 4)
$ws.Range("E2").Value = "mittelfristig"
$ws.Range("E3").Value = "mittelfristig"
$ws.Range("E4").Value = "kurzfristig"

# Then update Risiko (F) values to "kein Risiko" for rows 2-4
$ws.Range("F2").Value = "kein Risiko"
$ws.Range("F3").Value = "kein Risiko"
$ws.Range("F4").Value = "kein Risiko"

# Row 4 - Tagesgeld: Mindestanlagebetrag -> 0 (numeric)
$ws.Range("D4").Value = 0

# Move active selection to F4 (matches the saved selection in the diff)
$ws.Range("F4").Select()
